$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe the whole used range (values + formatting) so obsolete rows/styles
# (e.g. the old "strike-through" milestone rows) disappear, then rebuild
# the table from scratch with the new task list.
$ws.Range("B2:D20").Clear()

# Title
$ws.Range("B2").Value = "PHASE 1"
$ws.Range("B2").Font.Bold = $true

# Row 4
$ws.Range("B4").Value = "Test Bench "
$ws.Range("C4").Value = "Romain"

# Row 5
$ws.Range("B5").Value = "Thread pour les ennemis"
$ws.Range("D5").Value = "(Evt)"

# Row 6
$ws.Range("B6").Value = "Son"
$ws.Range("C6").Value = "Romain"

# Row 7
$ws.Range("B7").Value = "Graphisme / Animation"
$ws.Range("C7").Value = "Steve"

# Rows 8-19 (single column B)
$ws.Range("B8").Value = "Capacités spéciales (Decorator)"
$ws.Range("B9").Value = "Déplacement de blocs"
$ws.Range("B10").Value = "Interaction avec les acteurs (loutre)"
$ws.Range("B11").Value = "Création des ennemis et des niveaux"
$ws.Range("B12").Value = "Amélioration de l'histoire (Animations infos)"
$ws.Range("B13").Value = "Coffres et lieux secrets"
$ws.Range("B14").Value = "Trace de pas du pingouin dans la neige + Détection"
$ws.Range("B15").Value = "Déplacements intelligents des ennemis (IA)"
$ws.Range("B16").Value = "Champ de vision du pingouin perturbé (Brume, Luminosité)"
$ws.Range("B17").Value = "Niveau de difficulté haute"
$ws.Range("B18").Value = "Classement des meilleurs temps obtenus"
$ws.Range("B19").Value = "Générateur de maps pour les utilisateurs"

# Widen column B for the longer task descriptions (closest width the
# engine's char-width quantization allows to the authored 55.28515625).
$ws.Columns("B").ColumnWidth = 54.5

# Move the active selection to B20, matching the saved view state.
$ws.Range("B20").Select() | Out-Null
